$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SummaryReport")

$ws.Range("A6").Value = "C001878"
$ws.Range("B6").Value = "Brief - Appellant's Reply Brief"
$ws.Range("C6").Value = "Business Exception"
$ws.Range("D6").Value = "Document Processing Failure: Headings not found in the Document,Statement of Appealability,Statement of Facts. Case Number: C001878"
$ws.Range("E6").Value = "Failed"

$ws.Range("A7").Value = "C001878"
$ws.Range("B7").Value = "Brief - Appellant's Reply Brief"
$ws.Range("E7").Value = "Success"
